$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension to reflect the new used range (A1:T13)

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "F12"
$ws.Range("C2").Value = "Gp1ba"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1161036666666667
$ws.Range("H2").Value = 0.348311
$ws.Range("I2").Value = 0.4390792107664924
$ws.Range("J2").Value = 0.4390792107664924
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.6833285
$ws.Range("N2").Value = 7.366657
$ws.Range("O2").Value = 0.3319404283605227
$ws.Range("P2").Value = 0.2657978481314736
$ws.Range("Q2").Value = 0.4276479443878333
$ws.Range("R2").Value = 2.565887666327
$ws.Range("S2").Value = 0.1457481413060297
$ws.Range("T2").Value = 0.1167063093809994

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "F12"
$ws.Range("C3").Value = "Gp1ba"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1161036666666667
$ws.Range("H3").Value = 0.348311
$ws.Range("I3").Value = 0.4390792107664924
$ws.Range("J3").Value = 0.4390792107664924
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 1.145672333333333
$ws.Range("N3").Value = 3.437017
$ws.Range("O3").Value = 0.1032476373170262
$ws.Range("P3").Value = 0.1240117087834133
$ws.Range("Q3").Value = 0.1330167586985556
$ws.Range("R3").Value = 1.197150828287
$ws.Range("S3").Value = 0.04533389110666493
$ws.Range("T3").Value = 0.05445096321842519

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "F12"
$ws.Range("C4").Value = "Gp1ba"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1161036666666667
$ws.Range("H4").Value = 0.348311
$ws.Range("I4").Value = 0.4390792107664924
$ws.Range("J4").Value = 0.4390792107664924
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.5012456666666667
$ws.Range("N4").Value = 1.503737
$ws.Range("O4").Value = 0.04517210487937449
$ws.Range("P4").Value = 0.05425664025835297
$ws.Range("Q4").Value = 0.05819645980077778
$ws.Range("R4").Value = 0.523768138207
$ws.Range("S4").Value = 0.01983413215909698
$ws.Range("T4").Value = 0.02382296278347912

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "F12"
$ws.Range("C5").Value = "Gp1ba"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1161036666666667
$ws.Range("H5").Value = 0.348311
$ws.Range("I5").Value = 0.4390792107664924
$ws.Range("J5").Value = 0.4390792107664924
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.8904715
$ws.Range("N5").Value = 3.780943
$ws.Range("O5").Value = 0.1703687084965025
$ws.Range("P5").Value = 0.1364209726756327
$ws.Range("Q5").Value = 0.2194906728788333
$ws.Range("R5").Value = 1.316944037273
$ws.Range("S5").Value = 0.07480535806595094
$ws.Range("T5").Value = 0.05989961301441401

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "F12"
$ws.Range("C6").Value = "Gp1ba"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1161036666666667
$ws.Range("H6").Value = 0.348311
$ws.Range("I6").Value = 0.4390792107664924
$ws.Range("J6").Value = 0.4390792107664924
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.028466
$ws.Range("N6").Value = 9.085398
$ws.Range("O6").Value = 0.2729244218416247
$ws.Range("P6").Value = 0.3278120914029245
$ws.Range("Q6").Value = 0.3516160069753333
$ws.Range("R6").Value = 3.164544062777999
$ws.Range("S6").Value = 0.1198354397411218
$ws.Range("T6").Value = 0.1439354743729093

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "F12"
$ws.Range("C7").Value = "Gp1ba"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1161036666666667
$ws.Range("H7").Value = 0.348311
$ws.Range("I7").Value = 0.4390792107664924
$ws.Range("J7").Value = 0.4390792107664924
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8471700000000001
$ws.Range("N7").Value = 2.54151
$ws.Range("O7").Value = 0.07634669910494926
$ws.Range("P7").Value = 0.09170073874820307
$ws.Range("Q7").Value = 0.09835954329
$ws.Range("R7").Value = 0.88523588961
$ws.Range("S7").Value = 0.033522248387628
$ws.Range("T7").Value = 0.04026388799626531

# Row 8
$ws.Range("A8").Value = "Neutrophils"
$ws.Range("B8").Value = "F12"
$ws.Range("C8").Value = "Gp1ba"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1483216666666667
$ws.Range("H8").Value = 0.444965
$ws.Range("I8").Value = 0.5609207892335076
$ws.Range("J8").Value = 0.5609207892335076
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.6833285
$ws.Range("N8").Value = 7.366657
$ws.Range("O8").Value = 0.3319404283605227
$ws.Range("P8").Value = 0.2657978481314736
$ws.Range("Q8").Value = 0.5463174220008333
$ws.Range("R8").Value = 3.277904532005
$ws.Range("S8").Value = 0.186192287054493
$ws.Range("T8").Value = 0.1490915387504742

# Row 9
$ws.Range("A9").Value = "Neutrophils"
$ws.Range("B9").Value = "F12"
$ws.Range("C9").Value = "Gp1ba"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1483216666666667
$ws.Range("H9").Value = 0.444965
$ws.Range("I9").Value = 0.5609207892335076
$ws.Range("J9").Value = 0.5609207892335076
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 1.145672333333333
$ws.Range("N9").Value = 3.437017
$ws.Range("O9").Value = 0.1032476373170262
$ws.Range("P9").Value = 0.1240117087834133
$ws.Range("Q9").Value = 0.1699280299338889
$ws.Range("R9").Value = 1.529352269405
$ws.Range("S9").Value = 0.05791374621036131
$ws.Range("T9").Value = 0.06956074556498809

# Row 10
$ws.Range("A10").Value = "Neutrophils"
$ws.Range("B10").Value = "F12"
$ws.Range("C10").Value = "Gp1ba"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.1483216666666667
$ws.Range("H10").Value = 0.444965
$ws.Range("I10").Value = 0.5609207892335076
$ws.Range("J10").Value = 0.5609207892335076
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.5012456666666667
$ws.Range("N10").Value = 1.503737
$ws.Range("O10").Value = 0.04517210487937449
$ws.Range("P10").Value = 0.05425664025835297
$ws.Range("Q10").Value = 0.07434559268944445
$ws.Range("R10").Value = 0.669110334205
$ws.Range("S10").Value = 0.02533797272027752
$ws.Range("T10").Value = 0.03043367747487385

# Row 11
$ws.Range("A11").Value = "Neutrophils"
$ws.Range("B11").Value = "F12"
$ws.Range("C11").Value = "Gp1ba"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.1483216666666667
$ws.Range("H11").Value = 0.444965
$ws.Range("I11").Value = 0.5609207892335076
$ws.Range("J11").Value = 0.5609207892335076
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.8904715
$ws.Range("N11").Value = 3.780943
$ws.Range("O11").Value = 0.1703687084965025
$ws.Range("P11").Value = 0.1364209726756327
$ws.Range("Q11").Value = 0.2803978836658333
$ws.Range("R11").Value = 1.682387301995
$ws.Range("S11").Value = 0.0955633504305516
$ws.Range("T11").Value = 0.07652135966121865

# Row 12
$ws.Range("A12").Value = "Neutrophils"
$ws.Range("B12").Value = "F12"
$ws.Range("C12").Value = "Gp1ba"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.1483216666666667
$ws.Range("H12").Value = 0.444965
$ws.Range("I12").Value = 0.5609207892335076
$ws.Range("J12").Value = 0.5609207892335076
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.028466
$ws.Range("N12").Value = 9.085398
$ws.Range("O12").Value = 0.2729244218416247
$ws.Range("P12").Value = 0.3278120914029245
$ws.Range("Q12").Value = 0.4491871245633333
$ws.Range("R12").Value = 4.04268412107
$ws.Range("S12").Value = 0.1530889821005029
$ws.Range("T12").Value = 0.1838766170300151

# Row 13
$ws.Range("A13").Value = "Neutrophils"
$ws.Range("B13").Value = "F12"
$ws.Range("C13").Value = "Gp1ba"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.1483216666666667
$ws.Range("H13").Value = 0.444965
$ws.Range("I13").Value = 0.5609207892335076
$ws.Range("J13").Value = 0.5609207892335076
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.8471700000000001
$ws.Range("N13").Value = 2.54151
$ws.Range("O13").Value = 0.07634669910494926
$ws.Range("P13").Value = 0.09170073874820307
$ws.Range("Q13").Value = 0.12565366635
$ws.Range("R13").Value = 1.13088299715
$ws.Range("S13").Value = 0.04282445071732127
$ws.Range("T13").Value = 0.05143685075193776
